# Backup before dimension reduction:
# Decrement the numeric suffix of each "qN" label in column A (rows 2-97)
# so that q1 -> q0, q2 -> q1, ..., q96 -> q95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 97; $row++) {
    $n = $row - 1          # current numeric suffix (q1 at row2, q2 at row3, ...)
    $newValue = "q" + ($n - 1)
    $ws.Cells.Item($row, 1).Value = $newValue
}
